$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 35.125
$ws.Range("I5").Value = 35.125
$ws.Range("K5").Value = 35.125
$ws.Range("M5").Value = 79.875
$ws.Range("H11").Value = 87.666664
$ws.Range("I11").Value = 87.666664
$ws.Range("K11").Value = 87.666664
$ws.Range("M11").Value = 52.333336
$ws.Range("H33").Value = 1099.1428
$ws.Range("J33").Value = 3000
$ws.Range("L33").Value = 3000
$ws.Range("N33").Value = -3458
$ws.Range("H38").Value = 352.81818
$ws.Range("J38").Value = 925
$ws.Range("L38").Value = 2775
$ws.Range("N38").Value = -3519
$ws.Range("H41").Value = 1537.8334
$ws.Range("J41").Value = 1607
$ws.Range("L41").Value = 1607
$ws.Range("N41").Value = -2487
$ws.Range("H58").Value = 3466.6667
$ws.Range("J58").Value = 3466.6667
$ws.Range("L58").Value = 10400.0001
$ws.Range("N58").Value = -10700.0001
$ws.Range("H92").Value = 822.5
$ws.Range("I92").Value = 496.83334
$ws.Range("K92").Value = 496.83334
$ws.Range("M92").Value = 751.16666
$ws.Range("H116").Value = 6712.5713
$ws.Range("I116").Value = 5997.6
$ws.Range("K116").Value = 5997.6
$ws.Range("M116").Value = -2555.6
$ws.Range("H132").Value = 6096.524
$ws.Range("I132").Value = 5947.1665
$ws.Range("K132").Value = 17841.4995
$ws.Range("M132").Value = -15311.4995
$ws.Range("H133").Value = 185194
$ws.Range("J133").Value = 185194
$ws.Range("L133").Value = 185194
$ws.Range("N133").Value = -195314
$ws.Range("H137").Value = 2142.8096
$ws.Range("I137").Value = 1356.4667
$ws.Range("J137").Value = 4108.6665
$ws.Range("K137").Value = 4069.4001
$ws.Range("L137").Value = 12325.9995
$ws.Range("M137").Value = -1519.4001
$ws.Range("N137").Value = -17425.9995
$ws.Range("H138").Value = 3570.6956
$ws.Range("J138").Value = 3715.762
$ws.Range("L138").Value = 11147.286
$ws.Range("N138").Value = -21427.286

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 6473.3335
$ws.Range("I31").Value = 6473.3335
$ws.Range("K31").Value = 6473.3335
$ws.Range("M31").Value = -6179.3335
$ws.Range("H32").Value = 9340.272000000001
$ws.Range("I32").Value = 9340.272000000001
$ws.Range("K32").Value = 9340.272000000001
$ws.Range("M32").Value = -9053.272000000001
$ws.Range("H45").Value = 5000
$ws.Range("I45").Value = 5000
$ws.Range("K45").Value = 5000
$ws.Range("M45").Value = -4623
$ws.Range("H63").Value = 2625.5625
$ws.Range("I63").Value = 1612.1111
$ws.Range("J63").Value = 3928.5715
$ws.Range("K63").Value = 1612.1111
$ws.Range("L63").Value = 3928.5715
$ws.Range("M63").Value = -926.1111000000001
$ws.Range("N63").Value = -5300.5715
$ws.Range("H66").Value = 2625.5625
$ws.Range("I66").Value = 1612.1111
$ws.Range("J66").Value = 3928.5715
$ws.Range("K66").Value = 8060.5555
$ws.Range("L66").Value = 19642.8575
$ws.Range("M66").Value = -4628.5555
$ws.Range("N66").Value = -26506.8575
$ws.Range("H102").Value = 1568.1666
$ws.Range("I102").Value = 1568.1666
$ws.Range("K102").Value = 1568.1666
$ws.Range("M102").Value = 53.83339999999998
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 7062.4
$ws.Range("I132").Value = 7062.4
$ws.Range("K132").Value = 21187.2
$ws.Range("M132").Value = -18657.2
$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3481.923
$ws.Range("I105").Value = 2609.5715
$ws.Range("J105").Value = 4499.6665
$ws.Range("K105").Value = 2609.5715
$ws.Range("L105").Value = 4499.6665
$ws.Range("M105").Value = -862.5715
$ws.Range("N105").Value = -7993.6665

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2763.1365
$ws.Range("I31").Value = 2118.25
$ws.Range("J31").Value = 3537
$ws.Range("K31").Value = 2118.25
$ws.Range("L31").Value = 3537
$ws.Range("M31").Value = -1823.25
$ws.Range("N31").Value = -4127
$ws.Range("H34").Value = 2763.1365
$ws.Range("I34").Value = 2118.25
$ws.Range("J34").Value = 3537
$ws.Range("K34").Value = 2118.25
$ws.Range("L34").Value = 3537
$ws.Range("M34").Value = -1916.25
$ws.Range("N34").Value = -3941
$ws.Range("H44").Value = 15843.833
$ws.Range("I44").Value = 17765.75
$ws.Range("J44").Value = 12000
$ws.Range("K44").Value = 17765.75
$ws.Range("L44").Value = 12000
$ws.Range("M44").Value = -17323.75
$ws.Range("N44").Value = -12884
$ws.Range("H62").Value = 8800.666999999999
$ws.Range("I62").Value = 9012.333000000001
$ws.Range("K62").Value = 9012.333000000001
$ws.Range("M62").Value = -8388.333000000001
$ws.Range("H65").Value = 8800.666999999999
$ws.Range("I65").Value = 9012.333000000001
$ws.Range("K65").Value = 45061.665
$ws.Range("M65").Value = -41941.665
$ws.Range("H74").Value = 39999.168
$ws.Range("J74").Value = 39999.168
$ws.Range("L74").Value = 39999.168
$ws.Range("N74").Value = -41747.168
$ws.Range("H77").Value = 39999.168
$ws.Range("J77").Value = 39999.168
$ws.Range("L77").Value = 119997.504
$ws.Range("N77").Value = -128733.504
$ws.Range("H121").Value = 37499.5
$ws.Range("J121").Value = 35000
$ws.Range("L121").Value = 35000
$ws.Range("N121").Value = -37620

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1299.5714
$ws.Range("I11").Value = 275
$ws.Range("K11").Value = 825
$ws.Range("M11").Value = -685
$ws.Range("H134").Value = 2598.75
$ws.Range("I134").Value = 2598.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7796.25
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2726.25
$ws.Range("N134").ClearContents()
$ws.Range("H140").Value = 2273.9
$ws.Range("I140").Value = 1971
$ws.Range("K140").Value = 5913
$ws.Range("M140").Value = -733
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 10169.5
$ws.Range("I41").Value = 1525.5
$ws.Range("J41").Value = 11250
$ws.Range("K41").Value = 1525.5
$ws.Range("L41").Value = 11250
$ws.Range("M41").Value = -1170.5
$ws.Range("N41").Value = -11960
$ws.Range("H80").Value = 10000
$ws.Range("J80").Value = 10000
$ws.Range("L80").Value = 10000
$ws.Range("N80").Value = -11996
$ws.Range("H83").Value = 10000
$ws.Range("J83").Value = 10000
$ws.Range("L83").Value = 50000
$ws.Range("N83").Value = -59984
